$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" sheet (Worksheets.Item(1)): the existing row 2 ("2022-Q2") becomes
#    the new "2022-Q4" summary row, and a new row 3 is appended holding the
#    original "2022-Q2" summary row (same values/style the old row 2 had).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Row 3 = a formatted copy of row 2 (copies the bold/bordered A-column style),
# then filled in with the original "2022-Q2" totals.
$ws1.Range("A2:D2").Copy()
$ws1.Range("A3:D3").PasteSpecial(-4122)

$ws1.Cells.Item(3,1).Value = 1
$ws1.Cells.Item(3,2).Value = "2022-Q2"
$ws1.Cells.Item(3,3).Value = 1
$ws1.Cells.Item(3,4).Value = 0.21

# Row 2 is overwritten with the new "2022-Q4" totals.
$ws1.Cells.Item(2,2).Value = "2022-Q4"
$ws1.Cells.Item(2,4).Value = 0.2

# ---------------------------------------------------------------------------
# 2) Worksheet shuffle: the existing "2022-Q2" worksheet (Worksheets.Item(2))
#    is repurposed to hold the new "2022-Q4" detail data, and a brand new
#    worksheet named "2022-Q2" is inserted right after it holding the detail
#    data that used to live there. This keeps tab order 总计, 2022-Q4,
#    2022-Q2 and mirrors the sheetId/rId numbering of the target workbook.
# ---------------------------------------------------------------------------
$oldQ2Ws = $wb.Worksheets.Item(2)

$newQ2Ws = $wb.Worksheets.Add($null, $oldQ2Ws)
$newQ2Ws.Name = "Staging2022Q2"

# ---- populate the brand-new sheet with the original "2022-Q2" detail data
$newQ2Ws.Cells.Item(1,2).Value = "基金代码"
$newQ2Ws.Cells.Item(1,3).Value = "基金名称"
$newQ2Ws.Cells.Item(1,4).Value = "基金规模"
$newQ2Ws.Cells.Item(1,5).Value = "股票总仓位"
$newQ2Ws.Cells.Item(1,6).Value = "仓位占比"
$newQ2Ws.Cells.Item(1,7).Value = "持有市值(亿元)"
$newQ2Ws.Cells.Item(1,8).Value = "仓位排名"

$newQ2Ws.Cells.Item(2,1).Value = 0
$newQ2Ws.Cells.Item(2,2).Value = "'202801"
$newQ2Ws.Cells.Item(2,3).Value = "南方全球精选配置(QDII-FOF)"
$newQ2Ws.Cells.Item(2,4).Value = "'17.21"
$newQ2Ws.Cells.Item(2,5).Value = "'31.93"
$newQ2Ws.Cells.Item(2,6).Value = "'1.24"
$newQ2Ws.Cells.Item(2,7).Value = "'0.2134"
$newQ2Ws.Cells.Item(2,8).Value = 8

# strip the "quote prefix" styling that typing a leading apostrophe stamps on
# a cell, so the text cells keep the default (unstyled) look of the source.
$newQ2Ws.Range("B2").Style = "Normal"
$newQ2Ws.Range("D2:G2").Style = "Normal"

# header row + first data-row style (bold + border), matching the look the
# sheet had before the rename.
$oldQ2Ws.Range("B1:H1").Copy()
$newQ2Ws.Range("B1:H1").PasteSpecial(-4122)
$oldQ2Ws.Range("A2").Copy()
$newQ2Ws.Range("A2").PasteSpecial(-4122)

# ---- now overwrite the original sheet in place with the "2022-Q4" detail data
$oldQ2Ws.Cells.Item(2,2).Value = "'202801"
$oldQ2Ws.Cells.Item(2,3).Value = "南方全球精选配置（QDII-FOF）"
$oldQ2Ws.Cells.Item(2,4).Value = "'17.02"
$oldQ2Ws.Cells.Item(2,5).Value = "'32.64"
$oldQ2Ws.Cells.Item(2,6).Value = "'1.15"
$oldQ2Ws.Cells.Item(2,7).Value = "'0.1957"
$oldQ2Ws.Cells.Item(2,8).Value = 10

$oldQ2Ws.Range("B2").Style = "Normal"
$oldQ2Ws.Range("D2:G2").Style = "Normal"

# the "2022-Q4" sheet's header/first-data-row style is the bold+border look
# used on "总计" (style index 2), not the one the sheet had before (style
# index 1) - re-stamp it from ws1.
$ws1.Range("B1:D1").Copy()
$oldQ2Ws.Range("B1:H1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$oldQ2Ws.Range("A2").PasteSpecial(-4122)

# ---- finally rename the sheets into their target names.
$oldQ2Ws.Name = "2022-Q4"
$newQ2Ws.Name = "2022-Q2"
